$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H21").Value = 18999.75
$ws.Range("I21").Value = 0
$ws.Range("J21").Value = 18999.75
$ws.Range("K21").Value = 0
$ws.Range("L21").Value = 18999.75
$ws.Range("N21").Value = -19935.75
$ws.Range("M21").ClearContents()

$ws.Range("H23").Value = 18999.75
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 18999.75
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 18999.75
$ws.Range("N23").Value = -19467.75
$ws.Range("M23").ClearContents()

$ws.Range("H92").Value = 9260678
$ws.Range("I92").Value = 15874006
$ws.Range("J92").Value = 2019.2667
$ws.Range("K92").Value = 15874006
$ws.Range("L92").Value = 2019.2667
$ws.Range("M92").Value = -15872758
$ws.Range("N92").Value = -4515.2667

$ws.Range("H123").Value = 33000
$ws.Range("J123").Value = 33000
$ws.Range("L123").Value = 33000
$ws.Range("N123").Value = -42800

$ws.Range("H129").Value = 1447.1111
$ws.Range("J129").Value = 1890.84
$ws.Range("L129").Value = 5672.52
$ws.Range("N129").Value = -15672.52

$ws.Range("H137").Value = 1418.7234
$ws.Range("I137").Value = 1167.6129
$ws.Range("K137").Value = 3502.8387
$ws.Range("M137").Value = -952.8387000000002

$ws.Range("H140").Value = 70149.89999999999
$ws.Range("J140").Value = 70149.89999999999
$ws.Range("L140").Value = 70149.89999999999
$ws.Range("N140").Value = -80509.89999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 828.2
$ws.Range("I2").Value = 852.8570999999999
$ws.Range("J2").Value = 770.6667
$ws.Range("K2").Value = 852.8570999999999
$ws.Range("L2").Value = 770.6667
$ws.Range("M2").Value = -739.8570999999999
$ws.Range("N2").Value = -996.6667

$ws.Range("H32").Value = 5823678.5
$ws.Range("I32").Value = 6587115.5
$ws.Range("K32").Value = 6587115.5
$ws.Range("M32").Value = -6586828.5

$ws.Range("H46").Value = 4100
$ws.Range("J46").Value = 4100
$ws.Range("L46").Value = 4100
$ws.Range("N46").Value = -4738

$ws.Range("H63").Value = 35340.793
$ws.Range("I63").Value = 115969.75
$ws.Range("J63").Value = 4625
$ws.Range("K63").Value = 115969.75
$ws.Range("L63").Value = 4625
$ws.Range("M63").Value = -115283.75
$ws.Range("N63").Value = -5997

$ws.Range("H66").Value = 35340.793
$ws.Range("I66").Value = 115969.75
$ws.Range("J66").Value = 4625
$ws.Range("K66").Value = 579848.75
$ws.Range("L66").Value = 23125
$ws.Range("M66").Value = -576416.75
$ws.Range("N66").Value = -29989

$ws.Range("H74").Value = 3277.35
$ws.Range("I74").Value = 2203.182
$ws.Range("J74").Value = 4590.222
$ws.Range("K74").Value = 2203.182
$ws.Range("L74").Value = 4590.222
$ws.Range("M74").Value = -1329.182
$ws.Range("N74").Value = -6338.222

$ws.Range("H77").Value = 3277.35
$ws.Range("I77").Value = 2203.182
$ws.Range("J77").Value = 4590.222
$ws.Range("K77").Value = 11015.91
$ws.Range("L77").Value = 22951.11
$ws.Range("M77").Value = -6647.91
$ws.Range("N77").Value = -31687.11

$ws.Range("H116").Value = 828.2
$ws.Range("I116").Value = 852.8570999999999
$ws.Range("J116").Value = 770.6667
$ws.Range("K116").Value = 852.8570999999999
$ws.Range("L116").Value = 770.6667
$ws.Range("M116").Value = 1441.1429
$ws.Range("N116").Value = -5358.6667

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 828.2
$ws.Range("I3").Value = 852.8570999999999
$ws.Range("J3").Value = 770.6667
$ws.Range("K3").Value = 852.8570999999999
$ws.Range("L3").Value = 770.6667
$ws.Range("M3").Value = -738.8570999999999
$ws.Range("N3").Value = -998.6667

$ws.Range("H80").Value = 1838.8334
$ws.Range("J80").Value = 379.4
$ws.Range("L80").Value = 379.4
$ws.Range("N80").Value = -2375.4

$ws.Range("H83").Value = 1838.8334
$ws.Range("J83").Value = 379.4
$ws.Range("L83").Value = 1897
$ws.Range("N83").Value = -11881

$ws.Range("H134").Value = 2518.149
$ws.Range("I134").Value = 2357.838
$ws.Range("J134").Value = 3111.3
$ws.Range("K134").Value = 7073.514000000001
$ws.Range("L134").Value = 9333.900000000001
$ws.Range("M134").Value = -4538.514000000001
$ws.Range("N134").Value = -14403.9

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 888.5
$ws.Range("I87").Value = 888.5
$ws.Range("J87").Value = 0
$ws.Range("K87").Value = 2665.5
$ws.Range("L87").Value = 0
$ws.Range("M87").Value = -1417.5
$ws.Range("N87").ClearContents()

$ws.Range("H90").Value = 888.5
$ws.Range("I90").Value = 888.5
$ws.Range("J90").Value = 0
$ws.Range("K90").Value = 7996.5
$ws.Range("L90").Value = 0
$ws.Range("M90").Value = -1756.5
$ws.Range("N90").ClearContents()

$ws.Range("H131").Value = 3261.7917
$ws.Range("I131").Value = 406.47058
$ws.Range("J131").Value = 4827.613
$ws.Range("K131").Value = 1219.41174
$ws.Range("L131").Value = 14482.839
$ws.Range("M131").Value = 3820.58826
$ws.Range("N131").Value = -24562.839

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H68").Value = 61147.5
$ws.Range("I68").Value = 25000
$ws.Range("J68").Value = 97295
$ws.Range("K68").Value = 25000
$ws.Range("L68").Value = 97295
$ws.Range("M68").Value = -24189
$ws.Range("N68").Value = -98917

$ws.Range("H71").Value = 61147.5
$ws.Range("I71").Value = 25000
$ws.Range("J71").Value = 97295
$ws.Range("K71").Value = 75000
$ws.Range("L71").Value = 291885
$ws.Range("M71").Value = -70944
$ws.Range("N71").Value = -299997

$ws.Range("H102").Value = 2022.4
$ws.Range("I102").Value = 2056
$ws.Range("K102").Value = 2056
$ws.Range("M102").Value = -434

$ws.Range("H132").Value = 2293.6863
$ws.Range("I132").Value = 1720.1471
$ws.Range("J132").Value = 3440.7646
$ws.Range("K132").Value = 5160.4413
$ws.Range("L132").Value = 10322.2938
$ws.Range("M132").Value = -2630.4413
$ws.Range("N132").Value = -15382.2938

$ws.Range("H136").Value = 13039.634
$ws.Range("J136").Value = 10792.464
$ws.Range("L136").Value = 32377.392
$ws.Range("N136").Value = -37477.392

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H93").Value = 8256.823
$ws.Range("I93").Value = 9897.076999999999
$ws.Range("J93").Value = 2926
$ws.Range("K93").Value = 9897.076999999999
$ws.Range("L93").Value = 2926
$ws.Range("M93").Value = -8649.076999999999
$ws.Range("N93").Value = -5422

$ws.Range("H140").Value = 57192.285
$ws.Range("J140").Value = 57192.285
$ws.Range("L140").Value = 57192.285
$ws.Range("N140").Value = -67552.285

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H57").Value = 0
$ws.Range("J57").Value = 0
$ws.Range("L57").Value = 0
$ws.Range("N57").ClearContents()

$ws.Range("H81").Value = 3911.2666
$ws.Range("I81").Value = 3806.3635
$ws.Range("K81").Value = 7612.727
$ws.Range("M81").Value = -6551.727

$ws.Range("H84").Value = 3911.2666
$ws.Range("I84").Value = 3806.3635
$ws.Range("K84").Value = 38063.635
$ws.Range("M84").Value = -32759.635

$ws.Range("H96").Value = 4687.5
$ws.Range("I96").Value = 3957.1428
$ws.Range("J96").Value = 9800
$ws.Range("K96").Value = 3957.1428
$ws.Range("L96").Value = 9800
$ws.Range("M96").Value = -2584.1428
$ws.Range("N96").Value = -12546
